$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs, preserved formatting) ---
# A8: "Volume 32   Number  43" -> "...45" (last run, chars 21-22)
$ws.Range("A8").Characters(21, 2).Text = "45"
# C9: "Report Covering the Week  10/20/2025  Through  10/26/2025"
#     -> first date chars 27-36, second date (after length shrink) chars 47-56
$ws.Range("C9").Characters(27, 10).Text = "11/3/2025"
$ws.Range("C9").Characters(47, 10).Text = "11/9/2025"

# --- Data table updates (rows 15-31) ---
$ws.Range("F15").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("F15").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K15").Copy($ws.Range("E15"))
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 13
$ws.Range("L15").Value = 85.714285714285
$ws.Range("M15").Value = 116.666666666667
$ws.Range("N15").Value = -51.851851851851
$ws.Range("F15").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 54
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = -14.285714285714
$ws.Range("L16").Value = -10
$ws.Range("M16").Value = -40
$ws.Range("N16").Value = -89.090909090909
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 36.363636363636
$ws.Range("I17").Value = 125
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 17.924528301886
$ws.Range("L17").Value = 3.305785123966
$ws.Range("M17").Value = 37.362637362637
$ws.Range("N17").Value = -49.596774193548
$ws.Range("D18").Value = 1
$ws.Range("F18").Value = 2
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 63
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = -21.25
$ws.Range("L18").Value = -25.882352941176
$ws.Range("M18").Value = -70.967741935483
$ws.Range("N18").Value = -93.485005170630
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 333
$ws.Range("K19").Value = -2.402402402402
$ws.Range("L19").Value = -14.248021108179
$ws.Range("M19").Value = 16.487455197132
$ws.Range("N19").Value = -18.546365914787
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 109
$ws.Range("J20").Value = 169
$ws.Range("K20").Value = -35.502958579881
$ws.Range("L20").Value = -5.217391304347
$ws.Range("M20").Value = -16.153846153846
$ws.Range("N20").Value = -93.445580276608
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -25.675675675675
$ws.Range("I21").Value = 689
$ws.Range("J21").Value = 766
$ws.Range("K21").Value = -10.052219321148
$ws.Range("L21").Value = -10.403120936280
$ws.Range("M21").Value = -15.356265356265
$ws.Range("N21").Value = -81.882724165132
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = -22.222222222222
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -53.125
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = -35.714285714285
$ws.Range("I24").Value = 838
$ws.Range("J24").Value = 1225
$ws.Range("K24").Value = -31.591836734693
$ws.Range("L24").Value = -37.181409295352
$ws.Range("M24").Value = -19.731800766283
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 52
$ws.Range("H25").Value = -42.307692307692
$ws.Range("I25").Value = 391
$ws.Range("J25").Value = 751
$ws.Range("K25").Value = -47.936085219707
$ws.Range("L25").Value = -45.467224546722
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 30
$ws.Range("I26").Value = 347
$ws.Range("J26").Value = 310
$ws.Range("K26").Value = 11.935483870967
$ws.Range("L26").Value = 7.430340557275
$ws.Range("M26").Value = 0.872093023255
$ws.Range("F15").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("F15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 17.647058823529
$ws.Range("L27").Value = 66.666666666666
$ws.Range("C28").Value = 2
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 43
$ws.Range("K28").Value = -4.444444444444
$ws.Range("L28").Value = 65.384615384615
$ws.Range("C14").Copy($ws.Range("G31"))
$ws.Range("E14").Copy($ws.Range("H31"))
$ws.Range("L31").Value = -55.555555555555
